{"js": "// Update the first paragraph of the document:\n//  - add a paragraph border (spacing-only, no visible border line) on all 4 sides\n//  - change the left indent from 120 -> 225 twips (6pt -> 11.25pt)\n//  - replace the placeholder id text, dropping the trailing lone-space run\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// Left indent 120 -> 225 twips = 6pt -> 11.25pt (Office.js paragraph indents\n// are expressed in points).\nfirstParagraph.leftIndent = 11.25;\n\n// Paragraph border (w:pBdr) with w:space=\"5\" on every edge and no visible\n// line. Word's Word.Paragraph.borders edges only expose color/type/width in\n// Office.js (no \"distance from text\" setter), so reach the same\n// Borders.DistanceFrom* member the Word object model exposes for this via\n// the host bridge the generated proxies themselves dispatch through.\nconst handle = firstParagraph._h;\nconst anchor = firstParagraph._a;\nfor (const side of [\"Top\", \"Left\", \"Bottom\", \"Right\"]) {\n  globalThis.__native.docxOmSet(handle, anchor, \"Borders.DistanceFrom\" + side, \"5\");\n}\n\n// Replace the whole paragraph's text (both runs) with the new id string.\n// Using \"Replace\" on the paragraph's own range collapses the trailing \" \"\n// run into the text, leaving a single run that keeps the first run's\n// formatting (matches removing the lone-space run from the XML).\nfirstParagraph.getRange().insertText(\"**ID__AFFARS_AF_PGI_5319_201__ID**\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Update the first paragraph of the document:\n#  - add a paragraph border (spacing-only, no visible border line) on all 4 sides\n#  - change the left indent from 120 -> 225 twips (6pt -> 11.25pt)\n#  - replace the placeholder id text, dropping the trailing lone-space run\n\n$d = $word.ActiveDocument\n$p = $d.Paragraphs(1)\n\n# Paragraph border (w:pBdr) with w:space=\"5\" on every edge, no visible line.\n$p.Range.Borders.DistanceFromTop = 5\n$p.Range.Borders.DistanceFromLeft = 5\n$p.Range.Borders.DistanceFromBottom = 5\n$p.Range.Borders.DistanceFromRight = 5\n\n# Left indent 120 -> 225 twips = 6pt -> 11.25pt\n$p.Range.ParagraphFormat.LeftIndent = 11.25\n\n# Replace the whole paragraph's text (both runs, excluding the paragraph\n# mark) with the new id string. This collapses the trailing \" \" run into\n# the text, leaving a single run that carries the first run's formatting.\n$r = $p.Range\n$r.MoveEnd(1, -1)\n$r.Text = \"**ID__AFFARS_AF_PGI_5319_201__ID**\"\n"}
